# Generate Report for Handoff
#
# The localization-status report tracks, per source file / locale, when it
# was last handed off for translation. A new handoff run completed for
# "e0dedbf3-747d-4a26-9822-c10367072f15" (row 7 on every sheet), so its
# "Latest Handoff Date" / "Latest Handoff Datetime" columns are refreshed
# with the newer timestamp, both on the Overview roll-up sheet and on each
# per-locale detail sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column (D), row 7 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-26-11 16:26:14"

# --- zh-cn sheet: "Latest Handoff Datetime" column (E), row 7 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-11 16:26:11"

# --- de-de sheet: "Latest Handoff Datetime" column (E), row 7 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-11 16:26:14"
